$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Encode row (A23): now typechecks -> status OK, remove stale comment
$ws.Range("B23").Value = "OK"
$ws.Range("C23").ClearContents()

# StatefulLHAE row (A28): now typechecks -> status OK
$ws.Range("B28").Value = "OK"

# Move the active selection, matching the saved view state
$ws.Activate()
$ws.Range("C24").Select()
